$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F values for several rows
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4645
$wsExhibit.Range("F3").Value = 2526
$wsExhibit.Range("F4").Value = 480
$wsExhibit.Range("F10").Value = 191
$wsExhibit.Range("F12").Value = 1738
$wsExhibit.Range("F14").Value = 3873
$wsExhibit.Range("F15").Value = 37

# Sheet "全部类型" (sheet4): same logical rows but different row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4645
$wsAll.Range("F3").Value = 2526
$wsAll.Range("F4").Value = 480
$wsAll.Range("F12").Value = 191
$wsAll.Range("F16").Value = 1738
$wsAll.Range("F18").Value = 3873
$wsAll.Range("F19").Value = 37
